$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-25, columns B:N (G and M remain 0, unchanged)
$row2 = New-Object 'object[,]' 1,13
$row2[0,0] = 1.545629258784658
$row2[0,1] = 0.07717628820944356
$row2[0,2] = 0.01815585547026544
$row2[0,3] = 0.06961748551167801
$row2[0,4] = 3.570975794395338
$row2[0,5] = 0
$row2[0,6] = 0.07973214163530429
$row2[0,7] = 2.218182416072281
$row2[0,8] = 0.1477815485480001
$row2[0,9] = 1.352588289994941
$row2[0,10] = 0.3814899393340312
$row2[0,11] = 0
$row2[0,12] = 3.514942601425972
$ws.Range("B2:N2").Value = $row2

$row3 = New-Object 'object[,]' 1,13
$row3[0,0] = 1.51096712221613
$row3[0,1] = 0.07145143369959328
$row3[0,2] = 0.01780671062384087
$row3[0,3] = 0.06968500031635294
$row3[0,4] = 3.570356853167894
$row3[0,5] = 0
$row3[0,6] = 0.07973214163530429
$row3[0,7] = 2.221860651119997
$row3[0,8] = 0.1483247693326497
$row3[0,9] = 1.314252251397278
$row3[0,10] = 0.3783461650572946
$row3[0,11] = 0
$row3[0,12] = 3.534378708461375
$ws.Range("B3:N3").Value = $row3

$row4 = New-Object 'object[,]' 1,13
$row4[0,0] = 1.490490830023731
$row4[0,1] = 0.06797973033022231
$row4[0,2] = 0.01759941358128359
$row4[0,3] = 0.06974613591127365
$row4[0,4] = 3.571361696658954
$row4[0,5] = 0
$row4[0,6] = 0.07973214163530429
$row4[0,7] = 2.224912492693733
$row4[0,8] = 0.1486949130673629
$row4[0,9] = 1.291455622419903
$row4[0,10] = 0.3765853697630135
$row4[0,11] = 0
$row4[0,12] = 3.547150963245755
$ws.Range("B4:N4").Value = $row4

$row5 = New-Object 'object[,]' 1,13
$row5[0,0] = 1.482349752133899
$row5[0,1] = 0.06657583476555828
$row5[0,2] = 0.01751673526052855
$row5[0,3] = 0.06977601291391622
$row5[0,4] = 3.572119581773961
$row5[0,5] = 0
$row5[0,6] = 0.07973214163530429
$row5[0,7] = 2.226355728356623
$row5[0,8] = 0.1488549700800075
$row5[0,9] = 1.282352599688238
$row5[0,10] = 0.3759105295324048
$row5[0,11] = 0
$row5[0,12] = 3.552566584178663
$ws.Range("B5:N5").Value = $row5

$row6 = New-Object 'object[,]' 1,13
$row6[0,0] = 1.481010215716026
$row6[0,1] = 0.06634337275380631
$row6[0,2] = 0.01750311561298901
$row6[0,3] = 0.0697812742007109
$row6[0,4] = 3.572266475691478
$row6[0,5] = 0
$row6[0,6] = 0.07973214163530429
$row6[0,7] = 2.226607433454269
$row6[0,8] = 0.1488821047902888
$row6[0,9] = 1.280852338585191
$row6[0,10] = 0.3758010539817107
$row6[0,11] = 0
$row6[0,12] = 3.553478577166217
$ws.Range("B6:N6").Value = $row6

$row7 = New-Object 'object[,]' 1,13
$row7[0,0] = 1.490380213408883
$row7[0,1] = 0.06796075303866189
$row7[0,2] = 0.0175982912544832
$row7[0,3] = 0.06974651872499571
$row7[0,4] = 3.571370506880868
$row7[0,5] = 0
$row7[0,6] = 0.07973214163530429
$row7[0,7] = 2.22493114845971
$row7[0,8] = 0.1486970342978982
$row7[0,9] = 1.291332099246745
$row7[0,10] = 0.3765760956451913
$row7[0,11] = 0
$row7[0,12] = 3.547223146601048
$ws.Range("B7:N7").Value = $row7

$row8 = New-Object 'object[,]' 1,13
$row8[0,0] = 1.533510610707054
$row8[0,1] = 0.07519331872082091
$row8[0,2] = 0.01803401087269307
$row8[0,3] = 0.06963668851181204
$row8[0,4] = 3.570475042491921
$row8[0,5] = 0
$row8[0,6] = 0.07973214163530429
$row8[0,7] = 2.219286048031456
$row8[0,8] = 0.1479612634907639
$row8[0,9] = 1.339216181493555
$row8[0,10] = 0.3803708316214198
$row8[0,11] = 0
$row8[0,12] = 3.521470159716742
$ws.Range("B8:N8").Value = $row8

$row9 = New-Object 'object[,]' 1,13
$row9[0,0] = 1.624476395604859
$row9[0,1] = 0.08972391575396443
$row9[0,2] = 0.01894395548656291
$row9[0,3] = 0.06957686567080401
$row9[0,4] = 3.579703873641122
$row9[0,5] = 0
$row9[0,6] = 0.07973214163530429
$row9[0,7] = 2.214509288776497
$row9[0,8] = 0.1468081733530298
$row9[0,9] = 1.438999072898724
$row9[0,10] = 0.3891545921647008
$row9[0,11] = 0
$row9[0,12] = 3.477620897697619
$ws.Range("B9:N9").Value = $row9

$row10 = New-Object 'object[,]' 1,13
$row10[0,0] = 1.695197219961699
$row10[0,1] = 0.1006172201936408
$row10[0,2] = 0.01964551899122924
$row10[0,3] = 0.06962698538384515
$row10[0,4] = 3.593181707804874
$row10[0,5] = 0
$row10[0,6] = 0.07973214163530429
$row10[0,7] = 2.214836170887693
$row10[0,8] = 0.1461367586337623
$row10[0,9] = 1.51590036387887
$row10[0,10] = 0.3964240867529298
$row10[0,11] = 0
$row10[0,12] = 3.449459894757894
$ws.Range("B10:N10").Value = $row10

$row11 = New-Object 'object[,]' 1,13
$row11[0,0] = 1.728213758627703
$row11[0,1] = 0.1056214530942441
$row11[0,2] = 0.01997168734964916
$row11[0,3] = 0.06967005740891175
$row11[0,4] = 3.600768040855058
$row11[0,5] = 0
$row11[0,6] = 0.07973214163530429
$row11[0,7] = 2.215818099032816
$row11[0,8] = 0.1458692977196669
$row11[0,9] = 1.551666185535794
$row11[0,10] = 0.3999079873143216
$row11[0,11] = 0
$row11[0,12] = 3.437529257273908
$ws.Range("B11:N11").Value = $row11

$row12 = New-Object 'object[,]' 1,13
$row12[0,0] = 1.740837575793478
$row12[0,1] = 0.107523527283945
$row12[0,2] = 0.02009619368562809
$row12[0,3] = 0.06968926908410822
$row12[0,4] = 3.603850007860032
$row12[0,5] = 0
$row12[0,6] = 0.07973214163530429
$row12[0,7] = 2.216309729254121
$row12[0,8] = 0.1457734622181732
$row12[0,9] = 1.565322305470346
$row12[0,10] = 0.4012526427625573
$row12[0,11] = 0
$row12[0,12] = 3.433138004589424
$ws.Range("B12:N12").Value = $row12

$row13 = New-Object 'object[,]' 1,13
$row13[0,0] = 1.738113427806525
$row13[0,1] = 0.107113565732277
$row13[0,2] = 0.02006933509576569
$row13[0,3] = 0.06968500269331379
$row13[0,4] = 3.603176948491651
$row13[0,5] = 0
$row13[0,6] = 0.07973214163530429
$row13[0,7] = 2.216198520272542
$row13[0,8] = 0.1457938601075384
$row13[0,9] = 1.562376221060873
$row13[0,10] = 0.4009619195747405
$row13[0,11] = 0
$row13[0,12] = 3.434078106512658
$ws.Range("B13:N13").Value = $row13

$row14 = New-Object 'object[,]' 1,13
$row14[0,0] = 1.729249901385572
$row14[0,1] = 0.105777795645281
$row14[0,2] = 0.01998191074507361
$row14[0,3] = 0.06967157989918427
$row14[0,4] = 3.601017404084331
$row14[0,5] = 0
$row14[0,6] = 0.07973214163530429
$row14[0,7] = 2.215856144984599
$row14[0,8] = 0.145861304204935
$row14[0,9] = 1.552787431737272
$row14[0,10] = 0.4000181047162386
$row14[0,11] = 0
$row14[0,12] = 3.437165448096636
$ws.Range("B14:N14").Value = $row14

$row15 = New-Object 'object[,]' 1,13
$row15[0,0] = 1.723836504304757
$row15[0,1] = 0.1049605216955172
$row15[0,2] = 0.01992848969540262
$row15[0,3] = 0.06966373546567795
$row15[0,4] = 3.599721861198461
$row15[0,5] = 0
$row15[0,6] = 0.07973214163530429
$row15[0,7] = 2.215662030527753
$row15[0,8] = 0.145903324504669
$row15[0,9] = 1.546928650060039
$row15[0,10] = 0.3994432937436727
$row15[0,11] = 0
$row15[0,12] = 3.439073026349604
$ws.Range("B15:N15").Value = $row15

$row16 = New-Object 'object[,]' 1,13
$row16[0,0] = 1.693056488155662
$row16[0,1] = 0.1002911692215491
$row16[0,2] = 0.01962434303790772
$row16[0,3] = 0.06962457726601734
$row16[0,4] = 3.592715199444157
$row16[0,5] = 0
$row16[0,6] = 0.07973214163530429
$row16[0,7] = 2.214788764315387
$row16[0,8] = 0.1461550006884877
$row16[0,9] = 1.513578721114925
$row16[0,10] = 0.3961999607406312
$row16[0,11] = 0
$row16[0,12] = 3.450257303944866
$ws.Range("B16:N16").Value = $row16

$row17 = New-Object 'object[,]' 1,13
$row17[0,0] = 1.674390170281072
$row17[0,1] = 0.09743921800924227
$row17[0,2] = 0.01943954600876907
$row17[0,3] = 0.06960573728453312
$row17[0,4] = 3.588789500354778
$row17[0,5] = 0
$row17[0,6] = 0.07973214163530429
$row17[0,7] = 2.214466436963683
$row17[0,8] = 0.1463191108162825
$row17[0,9] = 1.493320016686027
$row17[0,10] = 0.3942555572567983
$row17[0,11] = 0
$row17[0,12] = 3.457343930555126
$ws.Range("B17:N17").Value = $row17

$row18 = New-Object 'object[,]' 1,13
$row18[0,0] = 1.663733382212968
$row18[0,1] = 0.09580344208046654
$row18[0,2] = 0.01933391771626347
$row18[0,3] = 0.0695968099772486
$row18[0,4] = 3.586668515986716
$row18[0,5] = 0
$row18[0,6] = 0.07973214163530429
$row18[0,7] = 2.214359471230587
$row18[0,8] = 0.1464170772582563
$row18[0,9] = 1.48174146194674
$row18[0,10] = 0.3931538493066853
$row18[0,11] = 0
$row18[0,12] = 3.461502789782728
$ws.Range("B18:N18").Value = $row18

$row19 = New-Object 'object[,]' 1,13
$row19[0,0] = 1.660138857206562
$row19[0,1] = 0.09525038312858669
$row19[0,2] = 0.01929826800973089
$row19[0,3] = 0.06959411563365592
$row19[0,4] = 3.585973914711545
$row19[0,5] = 0
$row19[0,6] = 0.07973214163530429
$row19[0,7] = 2.214336725353974
$row19[0,8] = 0.1464508613808739
$row19[0,9] = 1.477833829881035
$row19[0,10] = 0.3927836937087505
$row19[0,11] = 0
$row19[0,12] = 3.462925129776224
$ws.Range("B19:N19").Value = $row19

$row20 = New-Object 'object[,]' 1,13
$row20[0,0] = 1.676368997775512
$row20[0,1] = 0.0977423375003923
$row20[0,2] = 0.01945914956580808
$row20[0,3] = 0.06960754536445357
$row20[0,4] = 3.589193221939254
$row20[0,5] = 0
$row20[0,6] = 0.07973214163530429
$row20[0,7] = 2.214492632371318
$row20[0,8] = 0.1463012711518097
$row20[0,9] = 1.495468962614126
$row20[0,10] = 0.3944608184355047
$row20[0,11] = 0
$row20[0,12] = 3.456580974972226
$ws.Range("B20:N20").Value = $row20

$row21 = New-Object 'object[,]' 1,13
$row21[0,0] = 1.731850048422928
$row21[0,1] = 0.1061699510729568
$row21[0,2] = 0.02000756256430236
$row21[0,3] = 0.06967544387084956
$row21[0,4] = 3.601646037860405
$row21[0,5] = 0
$row21[0,6] = 0.07973214163530429
$row21[0,7] = 2.215953457853118
$row21[0,8] = 0.1458413465408945
$row21[0,9] = 1.555600843386884
$row21[0,10] = 0.4002946378256524
$row21[0,11] = 0
$row21[0,12] = 3.436255183875957
$ws.Range("B21:N21").Value = $row21

$row22 = New-Object 'object[,]' 1,13
$row22[0,0] = 1.768816223717977
$row22[0,1] = 0.1117191785282614
$row22[0,2] = 0.02037176525658069
$row22[0,3] = 0.06973672370039807
$row22[0,4] = 3.611003948694048
$row22[0,5] = 0
$row22[0,6] = 0.07973214163530429
$row22[0,7] = 2.217606442728723
$row22[0,8] = 0.1455724972204564
$row22[0,9] = 1.595555437930102
$row22[0,10] = 0.4042552658726066
$row22[0,11] = 0
$row22[0,12] = 3.42370919028086
$ws.Range("B22:N22").Value = $row22

$row23 = New-Object 'object[,]' 1,13
$row23[0,0] = 1.74902219628342
$row23[0,1] = 0.1087536528543751
$row23[0,2] = 0.02017685968394645
$row23[0,3] = 0.06970247533139329
$row23[0,4] = 3.605897913040664
$row23[0,5] = 0
$row23[0,6] = 0.07973214163530429
$row23[0,7] = 2.216660330012644
$row23[0,8] = 0.1457130876412158
$row23[0,9] = 1.574171064108924
$row23[0,10] = 0.4021278952443197
$row23[0,11] = 0
$row23[0,12] = 3.430337658783827
$ws.Range("B23:N23").Value = $row23

$row24 = New-Object 'object[,]' 1,13
$row24[0,0] = 1.675474137066828
$row24[0,1] = 0.09760528519478839
$row24[0,2] = 0.01945028488376366
$row24[0,3] = 0.06960672199987883
$row24[0,4] = 3.589010275904229
$row24[0,5] = 0
$row24[0,6] = 0.07973214163530429
$row24[0,7] = 2.214480545373888
$row24[0,8] = 0.1463093251956664
$row24[0,9] = 1.4944972109449
$row24[0,10] = 0.394367969529597
$row24[0,11] = 0
$row24[0,12] = 3.4569256434512
$ws.Range("B24:N24").Value = $row24

$row25 = New-Object 'object[,]' 1,13
$row25[0,0] = 1.599184536872514
$row25[0,1] = 0.08575519842273138
$row25[0,2] = 0.01869193384589352
$row25[0,3] = 0.06957647496794017
$row25[0,4] = 3.576031169168573
$row25[0,5] = 0
$row25[0,6] = 0.07973214163530429
$row25[0,7] = 2.21512782102392
$row25[0,8] = 0.1470891862452355
$row25[0,9] = 1.411374934514896
$row25[0,10] = 0.3866348840373774
$row25[0,11] = 0
$row25[0,12] = 3.488771166570785
$ws.Range("B25:N25").Value = $row25

